# Update the Timestamp column (D) so every row shares the same, newer
# run timestamp (a re-run of the login test suite).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-04-14 13:05:52"

$ws.Range("D2").Value = $newTimestamp
$ws.Range("D3").Value = $newTimestamp
$ws.Range("D4").Value = $newTimestamp
